# "Changed load cell termination pads"
#
# On the "Placement" sheet, three thin blank spacer rows (row height 6) are
# inserted to visually separate the "C1 / C2" pair, the "D1" row, and the
# "R1-R4" group from each other. Excel's native row Insert() shifts the
# existing rows (and their formulas/styles) down automatically, which is
# exactly the behaviour captured by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Placement")

# Insert a thin spacer row above the "C1" row (originally row 9).
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).RowHeight = 6

# "D1" (originally row 11) has now shifted down to row 12.
# Insert a thin spacer row above it.
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).RowHeight = 6

# "R1" (originally row 12) has now shifted down to row 14.
# Insert a thin spacer row above it.
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).RowHeight = 6

# Match the recorded cursor/selection position left behind in the file.
$ws.Activate()
$ws.Range("G24").Select()
